# Rename the sheet from "Property1" to "DataNode" (the concept unification
# mentioned in the commit message: DataNode / DataTable / Entity).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "DataNode"

# Update the active/frozen-pane selection to span the (now larger) editable
# area below the frozen header rows, A9:N35, instead of the single cell A9.
$ws.Range("A9:N35").Select() | Out-Null
